# Auto-generated Excel COM-interop script applying the Chocobo_Profits price/profit refresh
# across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 495.66666
$ws.Range("I2").Value = 598.2
$ws.Range("J2").Value = 367.5
$ws.Range("K2").Value = 598.2
$ws.Range("L2").Value = 367.5
$ws.Range("M2").Value = -485.2
$ws.Range("N2").Value = -593.5

$ws.Range("H4").Value = 145.2
$ws.Range("I4").Value = 116.888885
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 116.888885
$ws.Range("L4").Value = 400
$ws.Range("M4").Value = -2.888885000000002
$ws.Range("N4").Value = -628

$ws.Range("H18").Value = 252.14285
$ws.Range("I18").Value = 194.16667
$ws.Range("K18").Value = 194.16667
$ws.Range("M18").Value = 89.83332999999999

$ws.Range("H19").Value = 1481653.8
$ws.Range("J19").Value = 300
$ws.Range("L19").Value = 300
$ws.Range("N19").Value = -650

$ws.Range("H33").Value = 224.64102
$ws.Range("I33").Value = 190.85185
$ws.Range("J33").Value = 300.66666
$ws.Range("K33").Value = 190.85185
$ws.Range("L33").Value = 300.66666
$ws.Range("M33").Value = 38.14814999999999
$ws.Range("N33").Value = -758.66666

$ws.Range("H39").Value = 500.26666
$ws.Range("I39").Value = 200.44444
$ws.Range("J39").Value = 950
$ws.Range("K39").Value = 601.33332
$ws.Range("L39").Value = 2850
$ws.Range("M39").Value = -305.33332
$ws.Range("N39").Value = -3442

$ws.Range("H51").Value = 4998.3335
$ws.Range("I51").Value = 2000
$ws.Range("J51").Value = 5270.909
$ws.Range("K51").Value = 2000
$ws.Range("L51").Value = 5270.909
$ws.Range("M51").Value = -1516
$ws.Range("N51").Value = -6238.909

$ws.Range("H116").Value = 1010969
$ws.Range("I116").Value = 2501847.5
$ws.Range("J116").Value = 17050
$ws.Range("K116").Value = 2501847.5
$ws.Range("L116").Value = 17050
$ws.Range("M116").Value = -2498405.5
$ws.Range("N116").Value = -23934

$ws.Range("H125").Value = 2079.2942
$ws.Range("I125").Value = 2081.125
$ws.Range("J125").Value = 2077.6667
$ws.Range("K125").Value = 18730.125
$ws.Range("L125").Value = 18699.0003
$ws.Range("M125").Value = -16270.125
$ws.Range("N125").Value = -23619.0003

$ws.Range("H137").Value = 3425.3022
$ws.Range("I137").Value = 3379.9355
$ws.Range("J137").Value = 3542.5
$ws.Range("K137").Value = 10139.8065
$ws.Range("L137").Value = 10627.5
$ws.Range("M137").Value = -7589.806500000001
$ws.Range("N137").Value = -15727.5

$ws.Range("H138").Value = 2272.81
$ws.Range("I138").Value = 927.86664
$ws.Range("J138").Value = 2849.2144
$ws.Range("K138").Value = 2783.59992
$ws.Range("L138").Value = 8547.643199999999
$ws.Range("M138").Value = 2356.40008
$ws.Range("N138").Value = -18827.6432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6846.9487
$ws.Range("I32").Value = 4057.804
$ws.Range("K32").Value = 4057.804
$ws.Range("M32").Value = -3770.804

$ws.Range("H74").Value = 3019.102
$ws.Range("I74").Value = 2976.35
$ws.Range("J74").Value = 3209.111
$ws.Range("K74").Value = 2976.35
$ws.Range("L74").Value = 3209.111
$ws.Range("M74").Value = -2102.35
$ws.Range("N74").Value = -4957.111

$ws.Range("H77").Value = 3019.102
$ws.Range("I77").Value = 2976.35
$ws.Range("J77").Value = 3209.111
$ws.Range("K77").Value = 14881.75
$ws.Range("L77").Value = 16045.555
$ws.Range("M77").Value = -10513.75
$ws.Range("N77").Value = -24781.555

$ws.Range("H97").Value = 629.375
$ws.Range("I97").Value = 651.3333
$ws.Range("K97").Value = 651.3333
$ws.Range("M97").Value = -155.3333

$ws.Range("H132").Value = 2741.9333
$ws.Range("I132").Value = 1948.2258
$ws.Range("J132").Value = 4499.4287
$ws.Range("K132").Value = 5844.6774
$ws.Range("L132").Value = 13498.2861
$ws.Range("M132").Value = -3314.6774
$ws.Range("N132").Value = -18558.2861

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1062.875
$ws.Range("I94").Value = 1022.13043
$ws.Range("J94").Value = 2000
$ws.Range("K94").Value = 1022.13043
$ws.Range("L94").Value = 2000
$ws.Range("M94").Value = -571.13043
$ws.Range("N94").Value = -2902

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()

$ws.Range("H22").Value = 111722.836
$ws.Range("I22").Value = 166835.92
$ws.Range("K22").Value = 166835.92
$ws.Range("M22").Value = -166485.92

$ws.Range("H31").Value = 2881.7334
$ws.Range("I31").Value = 1187.92
$ws.Range("J31").Value = 4999
$ws.Range("K31").Value = 1187.92
$ws.Range("L31").Value = 4999
$ws.Range("M31").Value = -892.9200000000001
$ws.Range("N31").Value = -5589

$ws.Range("H34").Value = 2881.7334
$ws.Range("I34").Value = 1187.92
$ws.Range("J34").Value = 4999
$ws.Range("K34").Value = 1187.92
$ws.Range("L34").Value = 4999
$ws.Range("M34").Value = -985.9200000000001
$ws.Range("N34").Value = -5403

$ws.Range("H82").Value = 42400
$ws.Range("J82").Value = 42400
$ws.Range("L82").Value = 42400
$ws.Range("N82").Value = -43122

$ws.Range("H85").Value = 42400
$ws.Range("J85").Value = 42400
$ws.Range("L85").Value = 42400
$ws.Range("N85").Value = -44896

$ws.Range("H105").Value = 2664.4443
$ws.Range("I105").Value = 2568.5715
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 2568.5715
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -821.5715
$ws.Range("N105").Value = -6494

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H24").Value = 2650
$ws.Range("I24").Value = 300
$ws.Range("J24").Value = 5000
$ws.Range("K24").Value = 900
$ws.Range("L24").Value = 15000
$ws.Range("M24").Value = -670
$ws.Range("N24").Value = -15460

$ws.Range("H68").Value = 1460
$ws.Range("I68").Value = 1190
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 3570
$ws.Range("L68").Value = 6000
$ws.Range("M68").Value = -2759
$ws.Range("N68").Value = -7622

$ws.Range("H71").Value = 1460
$ws.Range("I71").Value = 1190
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 10710
$ws.Range("L71").Value = 18000
$ws.Range("M71").Value = -6654
$ws.Range("N71").Value = -26112

$ws.Range("H113").Value = 594
$ws.Range("I113").Value = 585.72974
$ws.Range("J113").Value = 617.53845
$ws.Range("K113").Value = 1757.18922
$ws.Range("L113").Value = 1852.61535
$ws.Range("M113").Value = 412.81078
$ws.Range("N113").Value = -6192.61535

$ws.Range("H129").Value = 3741.4443
$ws.Range("I129").Value = 4491.4287
$ws.Range("J129").Value = 1116.5
$ws.Range("K129").Value = 13474.2861
$ws.Range("L129").Value = 3349.5
$ws.Range("M129").Value = -8474.286100000001
$ws.Range("N129").Value = -13349.5

$ws.Range("H131").Value = 10643497
$ws.Range("J131").Value = 788.44183
$ws.Range("L131").Value = 2365.32549
$ws.Range("N131").Value = -12445.32549

$ws.Range("H132").Value = 2221.9656
$ws.Range("J132").Value = 4417.273
$ws.Range("L132").Value = 39755.457
$ws.Range("N132").Value = -44815.457

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6252.9766
$ws.Range("J70").Value = 8627.777
$ws.Range("L70").Value = 8627.777
$ws.Range("N70").Value = -9167.777

$ws.Range("H73").Value = 6252.9766
$ws.Range("J73").Value = 8627.777
$ws.Range("L73").Value = 8627.777
$ws.Range("N73").Value = -10499.777

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 2351
$ws.Range("J12").Value = 3702
$ws.Range("L12").Value = 3702
$ws.Range("N12").Value = -4042

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 14287512
$ws.Range("I81").Value = 17858390
$ws.Range("J81").Value = 4002
$ws.Range("K81").Value = 35716780
$ws.Range("L81").Value = 8004
$ws.Range("M81").Value = -35715719
$ws.Range("N81").Value = -10126

$ws.Range("H84").Value = 14287512
$ws.Range("I84").Value = 17858390
$ws.Range("J84").Value = 4002
$ws.Range("K84").Value = 178583900
$ws.Range("L84").Value = 40020
$ws.Range("M84").Value = -178578596
$ws.Range("N84").Value = -50628
